$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C25").Value = 1005
$ws.Range("D25").Value = 6009627
$ws.Range("E25").Value = 929.9948932219128
$ws.Range("G25").Value = 7.14285714285714
$ws.Range("H25").Value = 25.8601219808974
